$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---- Update "base" list (column F) ----
# Existing data occupies F2:F43 (42 entries). Read then re-write with the
# two new entries inserted in their alphabetically-sorted slot (right after
# "assertVarPresent(var)"), which pushes the remaining entries down by two
# rows (through F45).
$baseOld = @()
for ($r = 2; $r -le 43; $r++) {
    $baseOld += $ws.Cells.Item($r, 6).Value2
}

$baseNew = New-Object System.Collections.ArrayList
foreach ($item in $baseOld) {
    [void]$baseNew.Add($item)
    if ($item -eq "assertVarPresent(var)") {
        [void]$baseNew.Add("assertVarsNotPresent(vars)")
        [void]$baseNew.Add("assertVarsPresent(vars)")
    }
}

for ($i = 0; $i -lt $baseNew.Count; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $baseNew[$i]
}

# ---- Update "web" list (column AA) ----
# Existing data occupies AA2:AA151 (150 entries). Insert the new entry right
# after "assertLinkByLabel(label)", pushing the rest down by one row
# (through AA152).
$webOld = @()
for ($r = 2; $r -le 151; $r++) {
    $webOld += $ws.Cells.Item($r, 27).Value2
}

$webNew = New-Object System.Collections.ArrayList
foreach ($item in $webOld) {
    [void]$webNew.Add($item)
    if ($item -eq "assertLinkByLabel(label)") {
        [void]$webNew.Add("assertLocation(search)")
    }
}

for ($i = 0; $i -lt $webNew.Count; $i++) {
    $ws.Cells.Item($i + 2, 27).Value = $webNew[$i]
}

# ---- Update the defined names so the lists cover the new ranges ----
foreach ($n in $wb.Names) {
    if ($n.Name -eq "base") {
        $n.RefersTo = "='#system'!`$F`$2:`$F`$47"
    }
    if ($n.Name -eq "web") {
        $n.RefersTo = "='#system'!`$AA`$2:`$AA`$152"
    }
}

# Restore the originally-active sheet (editing the hidden "#system" sheet
# above leaves it as the active tab otherwise).
$wb.Worksheets.Item("MacroLibrary").Activate()
